$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the employee's name in row 3 ("Martin Zess" -> "Mark Zess")
$ws.Range("B3").Value = "Mark Zess"

# Update the active selection to A2
$ws.Range("A2").Select()
